# docs: correcoes no doc dicionario de dados
#
# Reworks the "Dicionario_de_Dados_Donate" workbook:
#  - lower-cases several table-title cells ("Tabela: ...")
#  - renames sheets to match the (lower_snake_case) table names
#  - adds three new sheets/tables: unidade_federativa, token_recuperacao, evento
#  - renames Agendamento -> doacao and corrects/extends its field list
#  - fixes a couple of wrong type/size values on the Notificacao sheet
#  - reorders the tabs to: municipio, unidade_federativa, token_recuperacao,
#    usuario, banco_de_leite, doacao, evento, Notificacao

$wb = $excel.ActiveWorkbook
$missing = [System.Reflection.Missing]::Value

function Set-Row([object]$ws, [int]$r, [object]$a, [object]$b, [object]$c, [object]$d) {
    $ws.Cells.Item($r, 1).Value = $a
    $ws.Cells.Item($r, 2).Value = $b
    $ws.Cells.Item($r, 3).Value = $c
    $ws.Cells.Item($r, 4).Value = $d
}

# ---------------------------------------------------------------------------
# 1) municipio (was "Municipio") - same sheet, just a title-case fix.
# ---------------------------------------------------------------------------
$shMunicipio = $wb.Worksheets.Item("Municipio")
$shMunicipio.Name = "municipio"
$shMunicipio.Range("A1").Value = "Tabela: município"

# ---------------------------------------------------------------------------
# 2) unidade_federativa - brand new sheet, cloned from municipio so it keeps
#    the same look & feel (fonts/borders/merged title row).
# ---------------------------------------------------------------------------
$shMunicipio.Copy($missing, $shMunicipio)
$shUf = $wb.Worksheets.Item($shMunicipio.Index + 1)
$shUf.Name = "unidade_federativa"

$shUf.Range("A1").Value = "Tabela: unidade_federativa"
Set-Row $shUf 3 "id"    "Serial"  "-"   "Chave primária da tabela municipio"
Set-Row $shUf 4 "nome"  "Varchar" 255   "Nome da Unidade federativa"
Set-Row $shUf 5 "sigla" "Varchar" 255   "Sigla Unidade federativa (Estado)"

$shUf.Columns.Item(4).ColumnWidth = 31.140625

# ---------------------------------------------------------------------------
# 3) token_recuperacao - brand new sheet, also cloned from municipio.
# ---------------------------------------------------------------------------
$shMunicipio.Copy($missing, $shUf)
$shToken = $wb.Worksheets.Item($shUf.Index + 1)
$shToken.Name = "token_recuperacao"

$shToken.Range("A1").Value = "Tabela: usuario_token"
Set-Row $shToken 3 "id"              "Serial"    "-" "Chave primária da tabela bancos_de_leite"
Set-Row $shToken 4 "usuario_id"      "Integer"   "-" "Usuario que esta solicitando a recuperação da senha"
Set-Row $shToken 5 "codigo"          "Varchar"   6   "código de recuperação"
Set-Row $shToken 6 "data_expiracao"  "Timestamp" "-" "data de expiração do token"
Set-Row $shToken 7 "usado"           "Boolean"   "-" "Se o token foi utilizado ou não"

$shToken.Columns.Item(4).ColumnWidth = 65

# ---------------------------------------------------------------------------
# 4) usuario (was "Usuario") - same sheet, correct a couple of field names.
# ---------------------------------------------------------------------------
$shUsuario = $wb.Worksheets.Item("Usuario")
$shUsuario.Name = "usuario"

# ---------------------------------------------------------------------------
# 5) banco_de_leite (was "Banco_Leite") - same sheet, title-case fix only.
# ---------------------------------------------------------------------------
$shBanco = $wb.Worksheets.Item("Banco_Leite")
$shBanco.Name = "banco_de_leite"
$shBanco.Range("A1").Value = "Tabela: banco_de_leite"

# ---------------------------------------------------------------------------
# 6) doacao (was "Agendamento") - same sheet, renamed + retitled; the field
#    list itself already matches what is needed for the "doacao" table.
# ---------------------------------------------------------------------------
$shDoacao = $wb.Worksheets.Item("Agendamento")
$shDoacao.Name = "doacao"
$shDoacao.Range("A1").Value = "Tabela: doacao"

# ---------------------------------------------------------------------------
# 7) evento - brand new sheet, cloned from doacao.
# ---------------------------------------------------------------------------
$shDoacao.Copy($missing, $shDoacao)
$shEvento = $wb.Worksheets.Item($shDoacao.Index + 1)
$shEvento.Name = "evento"

$shEvento.Range("A1").Value = "Tabela: evento"
Set-Row $shEvento 3 "id"           "Serial"  "-" "Chave primária da tabela agendamento."
Set-Row $shEvento 4 "titulo"       "Varchar" 255 "Titulo do evento."
Set-Row $shEvento 5 "descricao"    "Varchar" "Text" "descricao do evento"
Set-Row $shEvento 6 "data"         "Varchar" 255 "data do evento"
Set-Row $shEvento 7 "tipo"         "Varchar" 255 "Definição do tipo do evento se é noticia ou evento"
Set-Row $shEvento 8 "id_municipio" "Integer" "-" "Relacionamento com o municipio, para saber qual municipio vai ser o evento"

# row 9/10 (blank spacer row) from the doacao clone isn't part of "evento" -
# clear it out.
$shEvento.Rows.Item(9).Delete()
$shEvento.Rows.Item(9).Delete()

# ---------------------------------------------------------------------------
# 8) Notificacao - same sheet, stays last; fix two wrong type/size entries.
# ---------------------------------------------------------------------------
$shNotif = $wb.Worksheets.Item("Notificacao")
$shNotif.Cells.Item(4, 2).Value = "Boolean"
$shNotif.Cells.Item(4, 3).Value = 8
$shNotif.Cells.Item(5, 2).Value = "Int"
$shNotif.Cells.Item(5, 3).Value = 6

# Make sure Notificacao remains the last tab.
$shNotif.Move($missing, $shEvento)

# "token_recuperacao" is the active tab in the saved workbook - activate it
# last, since adding/copying sheets shifts the active tab around.
$shToken.Activate()
